$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.726.61'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.659.68'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.91'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.40'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.67%  '
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.399'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.86'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.13'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000194'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = '3.139.04'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").Value = '65.618.96'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").Value = '2.644.11'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '353.93'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.77'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +6.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.72'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '561.86'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.37%  '
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.69'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.422'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.57'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.98'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.50'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.71%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0619'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.36'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.35%  '
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.86'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("E50").Value = '  -7.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.816'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.23%  '
